$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-08-10 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-08-11 Sunday", 2)

# Update the multiplication table values. Each data row holds 5 equations,
# addressed by cell to avoid any ambiguity from duplicate values between cells.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="95×30="},
    @{Row=1;  Col=2; Text="15×76="},
    @{Row=1;  Col=3; Text="79×16="},
    @{Row=1;  Col=4; Text="74×56="},
    @{Row=1;  Col=5; Text="32×76="},

    @{Row=5;  Col=1; Text="77×85="},
    @{Row=5;  Col=2; Text="89×29="},
    @{Row=5;  Col=3; Text="55×30="},
    @{Row=5;  Col=4; Text="23×85="},
    @{Row=5;  Col=5; Text="66×75="},

    @{Row=10; Col=1; Text="52×35="},
    @{Row=10; Col=2; Text="40×56="},
    @{Row=10; Col=3; Text="56×33="},
    @{Row=10; Col=4; Text="29×41="},
    @{Row=10; Col=5; Text="79×80="},

    @{Row=15; Col=1; Text="43×72="},
    @{Row=15; Col=2; Text="42×39="},
    @{Row=15; Col=3; Text="67×29="},
    @{Row=15; Col=4; Text="46×92="},
    @{Row=15; Col=5; Text="16×26="},

    @{Row=20; Col=1; Text="46×16="},
    @{Row=20; Col=2; Text="87×39="},
    @{Row=20; Col=3; Text="83×44="},
    @{Row=20; Col=4; Text="14×19="},
    @{Row=20; Col=5; Text="18×42="}
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # equation text itself is replaced.
    $rng.MoveEnd(12, -1) | Out-Null
    $rng.Text = $u.Text
}
